# ActivityLogSheetWk4 - Week 4 activity log entries
# Fills in the activity-log table (name, week label, and four logged
# activities with dates/times/group hours) that was added to the
# previously-blank template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block: name + week label
$ws.Range("B2").Value2 = "Jesse Hare"
$ws.Range("F2").Value2 = "Week 4"

# Activity names (column A) first, so new shared strings land in the same
# order as the source: Complete Project Plan, then Team Meeting.
$ws.Range("A4").Value2 = "Complete Project Plan"
$ws.Range("A5").Value2 = "Team Meeting"
$ws.Range("A6").Value2 = "Complete Project Plan"
$ws.Range("A7").Value2 = "Complete Project Plan"

# Type (G/I) column
$ws.Range("C4").Value2 = "G"
$ws.Range("C5").Value2 = "G"
$ws.Range("C6").Value2 = "G"
$ws.Range("C7").Value2 = "G"

# Row 4: Complete Project Plan - Mon 19/08/2019, 9:00 AM - 2:00 PM, 5 group hrs
$ws.Range("D4").Value2 = 43696
$ws.Range("E4").Value2 = 0.375
$ws.Range("F4").Value2 = 0.58333333333333337
$ws.Range("G4").Value2 = 5

# Row 5: Team Meeting - Tue 20/08/2019, 9:00 AM - 11:00 AM, 2 group hrs
$ws.Range("D5").Value2 = 43697
$ws.Range("E5").Value2 = 0.375
$ws.Range("F5").Value2 = 0.45833333333333331
$ws.Range("G5").Value2 = 2

# Row 6: Complete Project Plan - Wed 21/08/2019, 9:00 AM - 3:00 PM, 6 individual hrs
$ws.Range("D6").Value2 = 43698
$ws.Range("E6").Value2 = 0.375
$ws.Range("F6").Value2 = 0.625
$ws.Range("H6").Value2 = 6

# Row 7: Complete Project Plan - Thu 22/08/2019, 9:00 AM - 4:00 PM, 7 individual hrs
$ws.Range("D7").Value2 = 43699
$ws.Range("E7").Value2 = 0.375
$ws.Range("F7").Value2 = 0.66666666666666663
$ws.Range("H7").Value2 = 7

# Leave the active selection on F7, matching the last-edited cell
$ws.Range("F7").Select()
